# Updated symbol list on Mon Dec 19 22:44:14 UTC 2022 with GitHub Actions
#
# Refreshes the coinranking.com crypto snapshot: rank #9 ("One"/ONE) jumped
# ahead of the #10-#17 exchange tokens (WazirX, MandalaExchangeToken, ...,
# CoinExToken), which all shifted down one rank, and every still-listed
# coin's Price cell (column D) was refreshed to the latest quote.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'241.00"
$ws.Range("D3").Value = "`'21.38"
$ws.Range("D4").Value = "`'5.135"
$ws.Range("D5").Value = "`'0.05550"
$ws.Range("D6").Value = "`'3.369"
$ws.Range("D7").Value = "`'6.351"
$ws.Range("D8").Value = "`'0.8050"
$ws.Range("D9").Value = "`'0.9475"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "`'0.01102"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "`'0.1393"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "`'0.07226"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "`'0.03068"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "`'0.03085"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "`'0.09281"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "`'3.616"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "`'0.001643"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "`'0.04706"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "`'0.006397"
$ws.Range("D20").Value = "`'0.004975"
$ws.Range("D21").Value = "`'0.001050"
$ws.Range("D22").Value = "`'0.0001505"
$ws.Range("D23").Value = "`'0.0003111"
$ws.Range("D24").Value = "`'3.750"
$ws.Range("D25").Value = "`'2.099"
$ws.Range("D26").Value = "`'0.3255"
$ws.Range("D27").Value = "`'0.1287"
$ws.Range("D40").Value = "`'0.03872"
$ws.Range("D41").Value = "`'0.006891"
$ws.Range("D42").Value = "`'0.1025"
$ws.Range("D43").Value = "`'0.003101"
$ws.Range("D44").Value = "`'0.008258"
$ws.Range("D45").Value = "`'0.00005961"
$ws.Range("D46").Value = "`'0.00000000753"
$ws.Range("D47").Value = "`'0.0005537"
$ws.Range("D48").Value = "`'0.6848"
$ws.Range("D49").Value = "`'0.1021"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").Value = "`'0.00002107"
$ws.Range("D51").Value = "`'0.01014"
